$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172 (shifts existing rows 172-229 down to 173-230),
# then populate it with the new weekly price entry.
$ws.Rows.Item(172).Insert()

$ws.Range("A172").Value = 7
$ws.Range("B172").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C172").Value = "Ñuble"
$ws.Range("D172").Value = 44559
$ws.Range("E172").Value = 16
$ws.Range("F172").Value = 100114013
$ws.Range("G172").Value = "Zanahoria"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 120
$ws.Range("K172").Value = 7500
$ws.Range("L172").Value = 8000
$ws.Range("M172").Value = 7750
$ws.Range("N172").Value = "$/saco 20 kilos"
$ws.Range("O172").Value = "Provincia de Diguillín"
$ws.Range("P172").Value = 388
$ws.Range("Q172").Value = 20
$ws.Range("R172").Value = "Hortaliza"
